$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IO_Map")

# New "Timer Unit" rows 11-15 appended below the existing I/O space map table.
$rows = @(
    @{ A = "0X08"; B = "W"; C = "us COUNT";  D = "D0-D7" },
    @{ A = "0X09"; B = "W"; C = "ms COUNT";  D = "D0-D7" },
    @{ A = "0X0A"; B = "W"; C = "sec COUNT"; D = "D0-D7" },
    @{ A = "0X0B"; B = "W"; C = "NOT USED";  D = "D0-D7" },
    @{ A = "0X08"; B = "R"; C = "STATUS";    D = "D0=RUNNING" }
)

$formatSource = $ws.Range("A7:D7")

$r = 11
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    # Match the highlighted formatting already used by the rows above (7-10).
    $formatSource.Copy()
    $destRow = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4))
    $destRow.PasteSpecial(-4122)

    $r = $r + 1
}

$excel.CutCopyMode = 0

# Widen column D to fit the longer "D0=RUNNING" bit label (stored OOXML
# width is ColumnWidth + ~0.83, so 12 here serializes to 12.83).
$ws.Columns.Item(4).ColumnWidth = 12

# Move the active selection as recorded after the edit.
$ws.Range("C5").Select() | Out-Null
